# Split the sole paragraph ("111222333") into four paragraphs:
#   1) "111222333"            (paragraph-mark formatting (rFonts hint) removed)
#   2) "333"
#   3) "444"
#   4) "555" (keeps the _GoBack bookmark that used to sit at the end of paragraph 1)
#
# We rebuild the first paragraph's range as raw WordprocessingML via InsertXML so we
# get exact control over paragraph marks, run formatting and bookmark placement.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = "<w:p $wNs>" +
           "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>111</w:t></w:r>" +
           "<w:r><w:t>222333</w:t></w:r>" +
       "</w:p>" +
       "<w:p $wNs><w:r><w:t>333</w:t></w:r></w:p>" +
       "<w:p $wNs><w:r><w:t>444</w:t></w:r></w:p>" +
       "<w:p $wNs>" +
           "<w:r><w:t>555</w:t></w:r>" +
           "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
           "<w:bookmarkEnd w:id=`"0`"/>" +
       "</w:p>"

$r.InsertXML($xml)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("Para" + $i + ": [" + $d.Paragraphs.Item($i).Range.Text + "]")
}
